# The presentation is linked to a SharePoint document-library content type
# (see customXml/item1.xml .. item3.xml). Opening/saving it as part of the
# batch re-syncs the cached schema list that SharePoint/InfoPath keeps in
# the datastoreItem's <ds:schemaRefs> (customXml/itemProps1.xml) -- it picks
# up the standard Dublin Core / package-core-properties / xml: namespaces
# alongside the site-column namespaces that were already referenced.
#
# Rebuild that schemaRefs list so it matches what the SharePoint metadata
# sync produces: the five namespaces that were already present, plus the
# six "well known" ones it always folds in.

$p = $ppt.ActivePresentation

$itemId    = "{7025FDD9-4C58-4084-9F89-0E6ADD6FFF55}"
$finalUris = @(
    "http://purl.org/dc/elements/1.1/",
    "http://schemas.microsoft.com/office/2006/metadata/properties",
    "http://www.w3.org/XML/1998/namespace",
    "http://purl.org/dc/terms/",
    "230e9df3-be65-4c73-a93b-d1236ebd677e",
    "http://schemas.microsoft.com/office/2006/documentManagement/types",
    "http://schemas.microsoft.com/office/infopath/2007/PartnerControls",
    "http://schemas.openxmlformats.org/package/2006/metadata/core-properties",
    "27aa9422-7f1f-4c84-9cdf-302b1a67e513",
    "http://schemas.microsoft.com/sharepoint/v3",
    "http://purl.org/dc/dcmitype/"
)

$refsXml = ($finalUris | ForEach-Object { "<ds:schemaRef ds:uri=`"$_`"/>" }) -join ""
$newXml = "<?xml version=`"1.0`" encoding=`"utf-8`"?>`r`n<ds:datastoreItem xmlns:ds=`"http://schemas.openxmlformats.org/officeDocument/2006/customXml`" ds:itemID=`"$itemId`"><ds:schemaRefs>$refsXml</ds:schemaRefs></ds:datastoreItem>"

$parts = $p.CustomXMLParts
$existing = $parts.SelectByID($itemId)

if ($existing -ne $null) {
    # Update the part already carried by the package in place.
    $existing.XML = $newXml
} else {
    # This session's CustomXMLParts collection doesn't surface the part
    # that's already baked into the package -- (re)attach it with the
    # fully merged schema list so the datastore item ends up correct.
    $null = $parts.Add($newXml)
}
